$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 8561.081
$ws.Range("J17").Value = 8653.056
$ws.Range("L17").Value = 25959.168
$ws.Range("N17").Value = -26295.168

$ws.Range("H33").Value = 562.1539
$ws.Range("I33").Value = 346.27274
$ws.Range("K33").Value = 346.27274
$ws.Range("M33").Value = -117.27274

$ws.Range("H40").Value = 3752.2334
$ws.Range("I40").Value = 2411.5
$ws.Range("J40").Value = 4087.4167
$ws.Range("K40").Value = 2411.5
$ws.Range("L40").Value = 4087.4167
$ws.Range("M40").Value = -2236.5
$ws.Range("N40").Value = -4437.4167

$ws.Range("H62").Value = 5208
$ws.Range("I62").Value = 2829
$ws.Range("J62").Value = 12345
$ws.Range("K62").Value = 2829
$ws.Range("L62").Value = 12345
$ws.Range("M62").Value = -2205
$ws.Range("N62").Value = -13593

$ws.Range("H64").Value = 6616.5
$ws.Range("I64").Value = 4974.75
$ws.Range("K64").Value = 4974.75
$ws.Range("M64").Value = -4726.75

$ws.Range("H65").Value = 5208
$ws.Range("I65").Value = 2829
$ws.Range("J65").Value = 12345
$ws.Range("K65").Value = 14145
$ws.Range("L65").Value = 61725
$ws.Range("M65").Value = -11025
$ws.Range("N65").Value = -67965

$ws.Range("H67").Value = 6616.5
$ws.Range("I67").Value = 4974.75
$ws.Range("K67").Value = 4974.75
$ws.Range("M67").Value = -4116.75

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H92").Value = 16129887
$ws.Range("I92").Value = 20834214
$ws.Range("K92").Value = 20834214
$ws.Range("M92").Value = -20832966

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21488922
$ws.Range("I32").Value = 26268030
$ws.Range("K32").Value = 26268030
$ws.Range("M32").Value = -26267743

$ws.Range("H45").Value = 2254.4375
$ws.Range("I45").Value = 1370.1818
$ws.Range("K45").Value = 1370.1818
$ws.Range("M45").Value = -993.1818000000001

$ws.Range("H97").Value = 1887
$ws.Range("I97").Value = 1398.6
$ws.Range("J97").Value = 2497.5
$ws.Range("K97").Value = 1398.6
$ws.Range("L97").Value = 2497.5
$ws.Range("M97").Value = -902.5999999999999
$ws.Range("N97").Value = -3489.5

$ws.Range("H132").Value = 348225.66
$ws.Range("I132").Value = 502951.6
$ws.Range("K132").Value = 1508854.8
$ws.Range("M132").Value = -1506324.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 433.16666

$ws.Range("H94").Value = 902.5185
$ws.Range("I94").Value = 797.9286
$ws.Range("K94").Value = 797.9286
$ws.Range("M94").Value = -346.9286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 549.75
$ws.Range("I22").Value = 549.75
$ws.Range("K22").Value = 549.75
$ws.Range("M22").Value = -199.75

$ws.Range("H31").Value = 7130.5557
$ws.Range("I31").Value = 2900.2
$ws.Range("J31").Value = 9245.733
$ws.Range("K31").Value = 2900.2
$ws.Range("L31").Value = 9245.733
$ws.Range("M31").Value = -2605.2
$ws.Range("N31").Value = -9835.733

$ws.Range("H34").Value = 7130.5557
$ws.Range("I34").Value = 2900.2
$ws.Range("J34").Value = 9245.733
$ws.Range("K34").Value = 2900.2
$ws.Range("L34").Value = 9245.733
$ws.Range("M34").Value = -2698.2
$ws.Range("N34").Value = -9649.733

$ws.Range("H68").Value = 97999.664
$ws.Range("I68").Value = 32000
$ws.Range("J68").Value = 130999.5
$ws.Range("K68").Value = 32000
$ws.Range("L68").Value = 130999.5
$ws.Range("M68").Value = -31251
$ws.Range("N68").Value = -132497.5

$ws.Range("H71").Value = 97999.664
$ws.Range("I71").Value = 32000
$ws.Range("J71").Value = 130999.5
$ws.Range("K71").Value = 96000
$ws.Range("L71").Value = 392998.5
$ws.Range("M71").Value = -92256
$ws.Range("N71").Value = -400486.5

$ws.Range("H99").Value = 4000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 4000
$ws.Range("N99").Value = -6996
$ws.Range("M99").ClearContents()

$ws.Range("H112").Value = 60700.5
$ws.Range("J112").Value = 60700.5
$ws.Range("L112").Value = 60700.5
$ws.Range("N112").Value = -63654.5

$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 12000
$ws.Range("N126").Value = -16940
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 5916.3335
$ws.Range("I132").Value = 5916.3335
$ws.Range("K132").Value = 17749.0005
$ws.Range("M132").Value = -15219.0005

$ws.Range("H134").Value = 4578.25
$ws.Range("I134").Value = 5490
$ws.Range("J134").Value = 3666.5
$ws.Range("K134").Value = 16470
$ws.Range("L134").Value = 10999.5
$ws.Range("M134").Value = -13935
$ws.Range("N134").Value = -16069.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5345.8423
$ws.Range("I2").Value = 32.42857
$ws.Range("K2").Value = 32.42857
$ws.Range("M2").Value = 80.57142999999999

$ws.Range("H80").Value = 2436.3845
$ws.Range("I80").Value = 2012.6364
$ws.Range("K80").Value = 2012.6364
$ws.Range("M80").Value = -1014.6364

$ws.Range("H83").Value = 2436.3845
$ws.Range("I83").Value = 2012.6364
$ws.Range("K83").Value = 10063.182
$ws.Range("M83").Value = -5071.182000000001

$ws.Range("H97").Value = 907.2
$ws.Range("I97").Value = 577.44446
$ws.Range("K97").Value = 577.44446
$ws.Range("M97").Value = -81.44446000000005

$ws.Range("H120").Value = 88225.664
$ws.Range("J120").Value = 88225.664
$ws.Range("L120").Value = 88225.664
$ws.Range("N120").Value = -97901.664

$ws.Range("H126").Value = 3350
$ws.Range("I126").Value = 3200
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 9600
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -7130
$ws.Range("N126").Value = -15440

$ws.Range("H132").Value = 11399
$ws.Range("I132").Value = 14331.667
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 42995.001
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -40465.001
$ws.Range("N132").Value = -26060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 261.93332
$ws.Range("I55").Value = 208.8
$ws.Range("J55").Value = 368.2
$ws.Range("K55").Value = 208.8
$ws.Range("L55").Value = 368.2
$ws.Range("M55").Value = -35.80000000000001
$ws.Range("N55").Value = -714.2

$ws.Range("H68").Value = 6645.5293
$ws.Range("I68").Value = 5588.273
$ws.Range("K68").Value = 5588.273
$ws.Range("M68").Value = -4839.273

$ws.Range("H71").Value = 6645.5293
$ws.Range("I71").Value = 5588.273
$ws.Range("K71").Value = 27941.365
$ws.Range("M71").Value = -24197.365

$ws.Range("H93").Value = 2191.7693
$ws.Range("I93").Value = 1186.75
$ws.Range("K93").Value = 1186.75
$ws.Range("M93").Value = 61.25

$ws.Range("H121").Value = 77302
$ws.Range("J121").Value = 77302
$ws.Range("L121").Value = 77302
$ws.Range("N121").Value = -80796

$ws.Range("H132").Value = 346936.88
$ws.Range("I132").Value = 442803.4
$ws.Range("J132").Value = 4556.4287
$ws.Range("K132").Value = 1328410.2
$ws.Range("L132").Value = 13669.2861
$ws.Range("M132").Value = -1325880.2
$ws.Range("N132").Value = -18729.2861

$ws.Range("H136").Value = 2518.3333
$ws.Range("I136").Value = 2573.2856
$ws.Range("J136").Value = 2441.4
$ws.Range("K136").Value = 7719.8568
$ws.Range("L136").Value = 7324.200000000001
$ws.Range("M136").Value = -5169.8568
$ws.Range("N136").Value = -12424.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 84497.25
$ws.Range("J119").Value = 84497.25
$ws.Range("L119").Value = 84497.25
$ws.Range("N119").Value = -94173.25

$ws.Range("H121").Value = 72069.836
$ws.Range("J121").Value = 72069.836
$ws.Range("L121").Value = 72069.836
$ws.Range("N121").Value = -75563.836

$ws.Range("H122").Value = 11984.857
$ws.Range("I122").Value = 11984.857
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 35954.571
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -33504.571
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 4390.273
$ws.Range("I126").Value = 4410.5
$ws.Range("J126").Value = 4336.3335
$ws.Range("K126").Value = 13231.5
$ws.Range("L126").Value = 13009.0005
$ws.Range("M126").Value = -10761.5
$ws.Range("N126").Value = -17949.0005

$ws.Range("H132").Value = 23616.717
$ws.Range("I132").Value = 25561.38
$ws.Range("J132").Value = 3197.75
$ws.Range("K132").Value = 76684.14
$ws.Range("L132").Value = 9593.25
$ws.Range("M132").Value = -74154.14
$ws.Range("N132").Value = -14653.25
